$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(21)
Write-Host "Before: $($shp.Name) left=$($shp.Left) top=$($shp.Top)"
$shp.Left = 350.93787401574804
$shp.Top = 125.36094488188976
Write-Host "After: $($shp.Name) left=$($shp.Left) top=$($shp.Top)"
